# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" between "2021-Q4" and "总计" holding the
#    per-fund breakdown for the new quarter.
# 2. Prepend a new row to the "总计" (totals) sheet summarising 2022-Q1
#    (4 funds, 3.83 billion RMB held), pushing the existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet.
# Clone "2021-Q4" (it already carries the exact header/style layout we need)
# and drop it right after "2021-Q4" / right before "总计".
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$newSheet = $wb.Worksheets.Item($q4.Index + 1)
$newSheet.Name = "2022-Q1"

# Row 2: 012348 / 天弘恒生科技指数型发起式证券投资基金（QDII）A
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "012348"
$newSheet.Range("C2").Value = "天弘恒生科技指数型发起式证券投资基金（QDII）A"
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "38.10"
$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "92.34"
$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "7.01"
$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "2.6708"
$newSheet.Range("H2").Value = 5

# Row 3: 012349 / 天弘恒生科技指数型发起式证券投资基金（QDII）C
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").NumberFormat = "@"
$newSheet.Range("B3").Value = "012349"
$newSheet.Range("C3").Value = "天弘恒生科技指数型发起式证券投资基金（QDII）C"
$newSheet.Range("D3").NumberFormat = "@"
$newSheet.Range("D3").Value = "14.77"
$newSheet.Range("E3").NumberFormat = "@"
$newSheet.Range("E3").Value = "92.34"
$newSheet.Range("F3").NumberFormat = "@"
$newSheet.Range("F3").Value = "7.01"
$newSheet.Range("G3").NumberFormat = "@"
$newSheet.Range("G3").Value = "1.0354"
$newSheet.Range("H3").Value = 5

# Extend the sheet with two more rows (copy-then-insert keeps the row-2/3
# formatting intact for the new rows).
$newSheet.Rows.Item(3).Copy()
$newSheet.Rows.Item(4).Insert()
$newSheet.Rows.Item(3).Copy()
$newSheet.Rows.Item(5).Insert()

# The row-insert above leaves the bordered "index" column (A) on a
# near-miss cloned style; re-stamp it from A3 (format only) so it matches
# the other index cells exactly.
$newSheet.Range("A3").Copy()
$newSheet.Range("A4").PasteSpecial(-4122)
$newSheet.Range("A5").PasteSpecial(-4122)

# Row 4: 002379 / 工银瑞信香港中小盘股票（QDII）人民币
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").NumberFormat = "@"
$newSheet.Range("B4").Value = "002379"
$newSheet.Range("C4").Value = "工银瑞信香港中小盘股票（QDII）人民币"
$newSheet.Range("D4").NumberFormat = "@"
$newSheet.Range("D4").Value = "1.84"
$newSheet.Range("E4").NumberFormat = "@"
$newSheet.Range("E4").Value = "86.48"
$newSheet.Range("F4").NumberFormat = "@"
$newSheet.Range("F4").Value = "3.45"
$newSheet.Range("G4").NumberFormat = "@"
$newSheet.Range("G4").Value = "0.0635"
$newSheet.Range("H4").Value = 10

# Row 5: 002380 / 工银瑞信香港中小盘股票（QDII）美元
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").NumberFormat = "@"
$newSheet.Range("B5").Value = "002380"
$newSheet.Range("C5").Value = "工银瑞信香港中小盘股票（QDII）美元"
$newSheet.Range("D5").NumberFormat = "@"
$newSheet.Range("D5").Value = "1.84"
$newSheet.Range("E5").NumberFormat = "@"
$newSheet.Range("E5").Value = "86.48"
$newSheet.Range("F5").NumberFormat = "@"
$newSheet.Range("F5").Value = "3.45"
$newSheet.Range("G5").NumberFormat = "@"
$newSheet.Range("G5").Value = "0.0635"
$newSheet.Range("H5").Value = 10

# ---------------------------------------------------------------------------
# Step 2: prepend a 2022-Q1 summary row to "总计", shifting the old rows down.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Stash the values currently sitting in row 2 (2021-Q4) before we move
# anything, so we can re-write them one row lower.
$oldRow2Date = $totalSheet.Range("B2").Value2
$oldRow2Count = $totalSheet.Range("C2").Value2
$oldRow2Value = $totalSheet.Range("D2").Value2

# Grow the table by one row at the bottom (keeps the index-column style
# consistent with the existing rows), then shuffle the data up/down.
$totalSheet.Rows.Item(3).Copy()
$totalSheet.Rows.Item(4).Insert()

# Row 4 now duplicates old row 3 (2021-Q3) - just fix its index cell. Also
# re-stamp its format (Insert leaves a near-miss cloned style on it) so it
# matches the other index cells exactly.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)
$totalSheet.Range("A4").Value = 2

# Row 3 becomes the old row 2 (2021-Q4).
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = $oldRow2Date
$totalSheet.Range("C3").Value = $oldRow2Count
$totalSheet.Range("D3").Value = $oldRow2Value

# Row 2 becomes the new 2022-Q1 summary.
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 3.83

Write-Output "2022-Q1 sheet added and 总计 updated"
